# Render website / remove theme (not needed) from docs:
# pull in the two extra Pandoc reference-doc paragraph styles
# ("Abstract Title", "Footnote Block Text") and retune Abstract's
# top spacing now that it follows a title line.

$d = $word.ActiveDocument
$wdStyleTypeParagraph = 1
$wdAlignParagraphCenter = 1

# 1) New "Abstract Title" style - centered, bold, small-caps-blue heading
#    that introduces the Abstract paragraph (Next -> Abstract).
$abstractTitle = $d.Styles.Add("AbstractTitle", $wdStyleTypeParagraph)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = $wdAlignParagraphCenter
$abstractTitle.ParagraphFormat.SpaceBefore = 15   ; # 300 twips
$abstractTitle.ParagraphFormat.SpaceAfter = 0     ; # 0 twips
$abstractTitle.Font.Size = 10                     ; # sz/szCs = 20 (half-points)
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 0x8A5A34              ; # BGR for w:color 345A8A

# 2) Abstract now follows Abstract Title, so it needs less space above it:
#    before 300 twips (15pt) -> 100 twips (5pt). "after" stays 300.
$d.Styles("Abstract").ParagraphFormat.SpaceBefore = 5   ; # 100 twips

# 3) New "Footnote Block Text" style - Block Text's indent/spacing treatment
#    applied to footnotes, based on (and followed by) Footnote Text.
$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", $wdStyleTypeParagraph)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = $d.Styles("FootnoteText")
$footnoteBlockText.NextParagraphStyle = $d.Styles("FootnoteText")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5    ; # 100 twips
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5     ; # 100 twips
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24    ; # 480 twips
$footnoteBlockText.ParagraphFormat.RightIndent = 24   ; # 480 twips
